$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.823.76'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.291.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '122.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.89%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.641'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.03%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("E11").Value = '  -0.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("E13").Value = '  +0.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.896'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.634.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.291.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.730.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.85%  '

$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.58'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.69%  '

$ws.Range("E25").Value = '  -4.65%  '

$ws.Range("E26").Value = '  +1.62%  '

$ws.Range("E27").Value = '  +1.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.32%  '

$ws.Range("E29").Value = '  +0.36%  '

$ws.Range("E30").Value = '  -0.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0916'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.52%  '

$ws.Range("E35").Value = '  +2.09%  '

$ws.Range("E36").Value = '  +4.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.69%  '

$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("E40").Value = '  +8.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '75.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.97%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.241'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.52%  '

$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("E46").Value = '  -0.92%  '

$ws.Range("E47").Value = '  +3.32%  '

$ws.Range("E48").Value = '  -2.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '73.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +36.88%  '

$ws.Range("E50").Value = '  +0.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '102.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.71%  '
